# Week 13 logging update for Cardinals Players Data workbook.
# Updates cumulative season totals on the "Rushing" and "Receiving" sheets.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$wsRush = $wb.Worksheets.Item("Rushing")

# Row 2 - K.Murray
$wsRush.Range("C2").Value = 24
$wsRush.Range("D2").Value = 20
$wsRush.Range("E2").Value = 11
$wsRush.Range("F2").Value = 22

# Row 5 - J.Conner
$wsRush.Range("C5").Value = 77
$wsRush.Range("D5").Value = 44
$wsRush.Range("E5").Value = 21
$wsRush.Range("F5").Value = 30

# Row 7 - E.Benjamin
$wsRush.Range("D7").Value = 14
$wsRush.Range("E7").Value = 20

# Row 9 - J.Ward
$wsRush.Range("C9").Value = 9
$wsRush.Range("D9").Value = 5
$wsRush.Range("F9").Value = 5

# --- Receiving sheet ---
$wsRecv = $wb.Worksheets.Item("Receiving")

# Row 3 - J.Conner
$wsRecv.Range("C3").Value = 15
$wsRecv.Range("D3").Value = 13
$wsRecv.Range("E3").Value = 3
$wsRecv.Range("F3").Value = 2

# Row 4 - C.Kirk
$wsRecv.Range("C4").Value = 3
$wsRecv.Range("D4").Value = 2

# Row 5 - D.Hopkins
$wsRecv.Range("C5").Value = 38
$wsRecv.Range("D5").Value = 30
$wsRecv.Range("E5").Value = 13
$wsRecv.Range("F5").Value = 7
$wsRecv.Range("G5").Value = 12
$wsRecv.Range("H5").Value = 8

# Row 6 - A.Green
$wsRecv.Range("C6").Value = 42
$wsRecv.Range("D6").Value = 25

# Row 7 - R.Moore
$wsRecv.Range("C7").Value = 43
$wsRecv.Range("D7").Value = 35
$wsRecv.Range("G7").Value = 6

# Row 8 - Z.Ertz
$wsRecv.Range("C8").Value = 48
$wsRecv.Range("D8").Value = 40

# Row 12 - D.Harris
$wsRecv.Range("C12").Value = 4
$wsRecv.Range("D12").Value = 2

# Row 13 - Z.Ertz
$wsRecv.Range("C13").Value = 37
$wsRecv.Range("D13").Value = 32

$wb.Save()
